$wb = $excel.ActiveWorkbook

# Rename the three "Acc_Upfront*" sheets so they rotate:
#   Acc_Upfront3 -> Acc_Upfront1
#   Acc_Upfront1 -> Acc_Upfront2
#   Acc_Upfront2 -> Acc_Upfront3
# Use a temporary name to avoid name collisions during the rename sequence.
$wsOldUpfront3 = $wb.Worksheets.Item("Acc_Upfront3")
$wsOldUpfront1 = $wb.Worksheets.Item("Acc_Upfront1")
$wsOldUpfront2 = $wb.Worksheets.Item("Acc_Upfront2")

$wsOldUpfront3.Name = "Acc_Upfront_tmp"
$wsOldUpfront1.Name = "Acc_Upfront1_new"
$wsOldUpfront2.Name = "Acc_Upfront2_new"

$wsOldUpfront3.Name = "Acc_Upfront1"
$wsOldUpfront1.Name = "Acc_Upfront2"
$wsOldUpfront2.Name = "Acc_Upfront3"

# Update the selected cell on the "Transactions" sheet.
$wsTransactions = $wb.Worksheets.Item("Transactions")
$wsTransactions.Range("D2").Select()

# Former "Acc_Upfront3" sheet (now named "Acc_Upfront1") loses its tab selection;
# its internal selection stays at G3.
$wsOldUpfront3.Range("G3").Select()

# Former "Acc_Upfront1" sheet (now named "Acc_Upfront2") selection moves to H22.
$wsOldUpfront1.Range("H22").Select()

# Former "Acc_Upfront2" sheet (now named "Acc_Upfront3") selection moves to I18,
# and this sheet becomes the active (selected) tab.
$wsOldUpfront2.Range("I18").Select()
$wsOldUpfront2.Activate()

$wb.Save()
